$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the per-week date headers with generic date placeholders,
# one placeholder text shared by all week columns for each service time.
$placeholder0700 = "{Day}, {dd} {MMM} {yyyy}`nPkl. 07.00 Wib,"
$placeholder1000 = "{Day}, {dd} {MMM} {yyyy}`nPkl. 10.00 Wib,"

$ws.Range("B2").Value = $placeholder0700
$ws.Range("C2").Value = $placeholder0700
$ws.Range("D2").Value = $placeholder0700
$ws.Range("E2").Value = $placeholder0700
$ws.Range("F2").Value = $placeholder0700

$ws.Range("B27").Value = $placeholder1000
$ws.Range("C27").Value = $placeholder1000
$ws.Range("D27").Value = $placeholder1000
$ws.Range("E27").Value = $placeholder1000
$ws.Range("F27").Value = $placeholder1000

# Update the active selection/view so the sheet opens scrolled to the top
# with G27 selected (instead of scrolled to row 10 with F28 selected).
$ws.Range("G27").Select()
